$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.661.56"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "1.883.35"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2938"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7399"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "1.883.15"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "30.646.15"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007566"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "2.130.96"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.340"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.248"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.919"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09734"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.294"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.170"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.57%  "
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7011"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01913"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.795"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.029"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8426"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.394"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.062"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "916.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  +2.00%  "
